$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 275, shifting rows 275:314 down to 276:315
$ws.Rows.Item(275).Insert()

# Populate the newly inserted row 275 with data
$ws.Cells.Item(275, 1).Value = 7
$ws.Cells.Item(275, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(275, 3).Value = "Ñuble"
$ws.Cells.Item(275, 4).Value = 45180
$ws.Cells.Item(275, 5).Value = 16
$ws.Cells.Item(275, 6).Value = 100112040
$ws.Cells.Item(275, 7).Value = "Cilantro"
$ws.Cells.Item(275, 8).Value = "Sin especificar"
$ws.Cells.Item(275, 9).Value = "Primera"
$ws.Cells.Item(275, 10).Value = 120
$ws.Cells.Item(275, 11).Value = 1300
$ws.Cells.Item(275, 12).Value = 1500
$ws.Cells.Item(275, 13).Value = 1400
$ws.Cells.Item(275, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(275, 15).Value = "Región de Ñuble"
$ws.Cells.Item(275, 16).Value = 1400
$ws.Cells.Item(275, 17).Value = 1
$ws.Cells.Item(275, 18).Value = "Hortaliza"
